# Generate Report for Handback
# Updates the localization-status workbook after a successful handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The de-de handback timestamp/version-mismatch errors are cleared since the
#    content is now back in sync, so the "Error Detail" column is blanked out
#    and the Latest Handback File / DateTime columns for de-de now point at the
#    freshly generated xlf.
#  - Column widths are adjusted: Status needs more room for the longer text,
#    Error Detail needs less room now that it is empty.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This string is shared across the Overview status columns (zh-cn/de-de) and the
# per-language "Status" column, so update every cell that currently shows it.
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: clear stale "version mismatch" error details ---
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# --- de-de sheet: refresh handback file/datetime, clear error details ---
$dede.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-13 14:43:52"
$dede.Range("P2").Value = ""

$dede.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-13 14:43:52"
$dede.Range("P3").Value = ""

# --- zh-cn: Latest Handback DateTime refresh ---
$zhcn.Range("K2").Value = "2016-08-13 14:43:43"
$zhcn.Range("K3").Value = "2016-08-13 14:43:43"

# --- Column width adjustments ---
# Overview: zh-cn/de-de status columns (E,F) need to be wider for the new text.
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# zh-cn / de-de: Status column (C) wider, Error Detail column (P) narrower.
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
